$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.007.68"
$ws.Range("E2").Value = "  +2.65%  "

$ws.Range("D3").Value = "2.342.85"
$ws.Range("E3").Value = "  +2.58%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'312.93"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("E6").Value = "  +3.31%  "

$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "'0.621"
$ws.Range("E9").Value = "  +3.32%  "

$ws.Range("D10").Value = "'41.36"
$ws.Range("E10").Value = "  +4.74%  "

$ws.Range("D11").Value = "'0.0920"
$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("E12").Value = "  +2.79%  "

$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("E14").Value = "  +1.89%  "

$ws.Range("D15").Value = "'15.54"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").Value = "2.697.59"

$ws.Range("D17").Value = "2.339.03"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("D18").Value = "43.944.55"
$ws.Range("E18").Value = "  +2.76%  "

$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = "  +3.17%  "

$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("D21").Value = "'12.95"
$ws.Range("E21").Value = "  -5.50%  "

$ws.Range("D22").Value = "'74.40"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "'3.48"
$ws.Range("E23").Value = "  -1.04%  "

$ws.Range("D24").Value = "'269.10"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  +4.22%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'7.53"
$ws.Range("E27").Value = "  +5.70%  "

$ws.Range("D28").Value = "'11.18"
$ws.Range("E28").Value = "  +3.48%  "

$ws.Range("E29").Value = "  -1.87%  "

$ws.Range("D30").Value = "'39.11"
$ws.Range("E30").Value = "  +6.00%  "

$ws.Range("D31").Value = "'22.70"
$ws.Range("E31").Value = "  +1.20%  "

$ws.Range("D32").Value = "'168.58"
$ws.Range("E32").Value = "  +1.37%  "

$ws.Range("D33").Value = "'0.0891"
$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("E34").Value = "  +8.18%  "

$ws.Range("E35").Value = "  +1.52%  "

$ws.Range("D36").Value = "'4.77"
$ws.Range("E36").Value = "  +5.56%  "

$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("E38").Value = "  +4.28%  "

$ws.Range("D39").Value = "'2.90"
$ws.Range("E39").Value = "  +9.19%  "

$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("D41").Value = "'1.72"
$ws.Range("E41").Value = "  +9.27%  "

$ws.Range("D42").Value = "'104.74"
$ws.Range("E42").Value = "  +10.68%  "

$ws.Range("D43").Value = "'0.239"
$ws.Range("E43").Value = "  +3.96%  "

$ws.Range("D44").Value = "'71.85"
$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").Value = "'13.28"
$ws.Range("E45").Value = "  +9.55%  "

$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").Value = "'114.61"
$ws.Range("E47").Value = "  +1.20%  "

$ws.Range("D48").Value = "1.663.47"
$ws.Range("E48").Value = "  -3.99%  "

$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'76.95"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'9.00"
$ws.Range("E50").Value = "  +3.32%  "

$ws.Range("B51").Value = "MinaProtocolToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D51").Value = "'1.56"
$ws.Range("E51").Value = "  +8.43%  "
